$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.225.35"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.42%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "1.884.53"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.01%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'244.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.32%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'0.688"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.84%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'42.51"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.84%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("E9").Value = "  +3.44%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'54.96"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.89%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.0739"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.16%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'0.0981"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.53%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("E13").Value = "  +8.18%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'0.774"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +10.02%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "2.158.19"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.82%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("E16").Value = "  +3.18%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "1.883.45"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.99%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "35.212.04"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.44%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'73.05"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.56%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "0.0₃0820"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.03%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'243.20"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.93%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("E22").Value = "  +2.26%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'5.13"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.25%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'2.64"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.46%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.15%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("B26").Value = "Monero"
$ws.Range("B26").Style = "Normal"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("C26").Style = "Normal"
$ws.Range("D26").Value = "'166.89"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.51%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("B27").Style = "Normal"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("C27").Style = "Normal"
$ws.Range("D27").Value = "'2.13"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.74%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'8.47"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.92%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'18.20"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.33%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("E30").Value = "  +1.03%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("B31").Value = "Filecoin"
$ws.Range("B31").Style = "Normal"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("C31").Style = "Normal"
$ws.Range("D31").Value = "'4.26"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.81%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("B32").Value = "Hedera"
$ws.Range("B32").Style = "Normal"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("C32").Style = "Normal"
$ws.Range("D32").Value = "'0.0591"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.77%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'1.87"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +19.42%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'4.15"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.01%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("E36").Value = "  -13.09%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("E37").Value = "  +3.96%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("E38").Value = "  -0.78%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'0.0718"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.70%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'0.0220"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.45%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'97.59"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.25%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'17.05"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.58%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("E43").Value = "  +0.76%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "1.321.66"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.01%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'13.17"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +11.97%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("E46").Value = "  +2.91%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'0.0809"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'2.73"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.91%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'6.24"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.09%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "2.057.60"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.42%  "
$ws.Range("E51").Style = "Normal"
